# Insert a new data row at row 251 (pushes the existing rows 251-364 down
# to 252-365, growing the used range from A1:R364 to A1:R365), then fill
# the newly inserted row with its own record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("251:251").Insert()

$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 44726
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112032
$ws.Range("G251").Value = "Zapallo italiano"
$ws.Range("H251").Value = "Sin especificar"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 300
$ws.Range("K251").Value = 12000
$ws.Range("L251").Value = 12000
$ws.Range("M251").Value = 12000
$ws.Range("N251").Value = '$/caja 50 unidades'
$ws.Range("O251").Value = "Región de Arica y Parinacota"
$ws.Range("P251").Value = 240
$ws.Range("Q251").Value = 50
$ws.Range("R251").Value = "Hortaliza"
